$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume figures.
# Row 46 (BabyDogeCoin) was delisted, so rows 47-51 shift up one slot
# and a new "Cronos" entry is appended as the new row 51.

# Row 2
$ws.Range('D2').Value = '29.401.23'
$ws.Range('E2').Value = '  +0.43%  '

# Row 3
$ws.Range('D3').Value = '1.843.73'
$ws.Range('E3').Value = '  +0.22%  '

# Row 4
$ws.Range('D4').Value = "'0.9990"
$ws.Range('E4').Value = '  +0.11%  '

# Row 5
$ws.Range('D5').Value = "'240.25"
$ws.Range('E5').Value = '  +0.16%  '

# Row 6
$ws.Range('D6').Value = "'0.6341"
$ws.Range('E6').Value = '  +1.37%  '

# Row 7
$ws.Range('D7').Value = "'0.9999"
$ws.Range('E7').Value = '  +0.07%  '

# Row 8
$ws.Range('D8').Value = "'0.07476"
$ws.Range('E8').Value = '  -0.04%  '

# Row 9
$ws.Range('D9').Value = "'25.14"
$ws.Range('E9').Value = '  +3.79%  '

# Row 10
$ws.Range('D10').Value = "'0.2906"
$ws.Range('E10').Value = '  +0.42%  '

# Row 11
$ws.Range('E11').Value = '  +0.44%  '

# Row 12
$ws.Range('D12').Value = '1.872.19'
$ws.Range('E12').Value = '  +1.79%  '

# Row 13
$ws.Range('D13').Value = "'4.987"
$ws.Range('E13').Value = '  +0.12%  '

# Row 14
$ws.Range('D14').Value = "'0.6790"
$ws.Range('E14').Value = '  +0.32%  '

# Row 15
$ws.Range('D15').Value = "'0.00001021"
$ws.Range('E15').Value = '  -0.56%  '

# Row 16
$ws.Range('D16').Value = "'82.02"
$ws.Range('E16').Value = '  -0.01%  '

# Row 17
$ws.Range('D17').Value = "'6.251"
$ws.Range('E17').Value = '  +2.70%  '

# Row 18
$ws.Range('D18').Value = '29.463.51'
$ws.Range('E18').Value = '  +0.47%  '

# Row 19
$ws.Range('D19').Value = "'230.06"
$ws.Range('E19').Value = '  +0.70%  '

# Row 20
$ws.Range('E20').Value = '  +0.80%  '

# Row 21
$ws.Range('E21').Value = '  +0.06%  '

# Row 22
$ws.Range('D22').Value = "'7.408"
$ws.Range('E22').Value = '  +0.48%  '

# Row 23
$ws.Range('D23').Value = "'1.001"
$ws.Range('E23').Value = '  +0.16%  '

# Row 24
$ws.Range('D24').Value = "'158.24"
$ws.Range('E24').Value = '  +0.17%  '

# Row 25
$ws.Range('D25').Value = "'8.507"
$ws.Range('E25').Value = '  +1.69%  '

# Row 26
$ws.Range('E26').Value = '  -1.27%  '

# Row 27
$ws.Range('E27').Value = '  -0.16%  '

# Row 28
$ws.Range('D28').Value = "'0.06552"
$ws.Range('E28').Value = '  +15.14%  '

# Row 29
$ws.Range('D29').Value = "'1.437"
$ws.Range('E29').Value = '  +2.68%  '

# Row 30
$ws.Range('D30').Value = "'1.490"
$ws.Range('E30').Value = '  +1.39%  '

# Row 31
$ws.Range('D31').Value = "'4.072"
$ws.Range('E31').Value = '  -0.55%  '

# Row 32
$ws.Range('D32').Value = "'4.050"
$ws.Range('E32').Value = '  +0.54%  '

# Row 33
$ws.Range('D33').Value = "'1.843"
$ws.Range('E33').Value = '  +1.54%  '

# Row 34
$ws.Range('E34').Value = '  +0.01%  '

# Row 35
$ws.Range('D35').Value = "'0.7012"
$ws.Range('E35').Value = '  +1.49%  '

# Row 36
$ws.Range('D36').Value = "'2.579"
$ws.Range('E36').Value = '  -0.07%  '

# Row 37
$ws.Range('D37').Value = "'0.01859"
$ws.Range('E37').Value = '  +2.88%  '

# Row 38
$ws.Range('D38').Value = '1.249.74'
$ws.Range('E38').Value = '  +0.57%  '

# Row 39
$ws.Range('D39').Value = "'2.818"
$ws.Range('E39').Value = '  +0.19%  '

# Row 40
$ws.Range('D40').Value = "'6.751"

# Row 41
$ws.Range('D41').Value = "'0.9423"
$ws.Range('E41').Value = '  +4.20%  '

# Row 42
$ws.Range('D42').Value = "'0.9997"
$ws.Range('E42').Value = '  +0.21%  '

# Row 43
$ws.Range('D43').Value = '2.010.49'
$ws.Range('E43').Value = '  +0.51%  '

# Row 44
$ws.Range('D44').Value = "'101.23"
$ws.Range('E44').Value = '  -0.11%  '

# Row 45
$ws.Range('E45').Value = '  -0.20%  '

# Row 46
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = "'7.067"
$ws.Range('E46').Value = '  +0.17%  '

# Row 47
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = "'1.718"
$ws.Range('E47').Value = '  +3.92%  '

# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'9.032"
$ws.Range('E48').Value = '  +0.63%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = "'0.1148"
$ws.Range('E49').Value = '  -1.58%  '

# Row 50
$ws.Range('B50').Value = 'TheSandbox'
$ws.Range('C50').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D50').Value = "'0.3921"
$ws.Range('E50').Value = '  -0.35%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = "'0.05679"
$ws.Range('E51').Value = '  -0.15%  '
